$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shorten the Scene FilePath values from "../../NFDataCfg/..." to "../NFDataCfg/..."
# (the project's Excel_Ini folder moved one level, so the relative path needs one
# fewer "../" to reach NFDataCfg) for all six scene rows.
$ws.Range("F10").Value = "../NFDataCfg/Ini/Scene/1.xml"
$ws.Range("F11").Value = "../NFDataCfg/Ini/Scene/2.xml"
$ws.Range("F12").Value = "../NFDataCfg/Ini/Scene/3.xml"
$ws.Range("F13").Value = "../NFDataCfg/Ini/Scene/4.xml"
$ws.Range("F14").Value = "../NFDataCfg/Ini/Scene/5.xml"
$ws.Range("F15").Value = "../NFDataCfg/Ini/Scene/6.xml"

# Widen the FilePath column so the (now differently sized) paths are fully visible.
$ws.Columns.Item(6).ColumnWidth = 31.22

# Leave the cursor where the author last left it when saving.
$ws.Range("F23").Select() | Out-Null
